# Updated cryptos list on Sat Nov 30 18:52:52 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for each coin row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '96.954.35'
$ws.Range("E2").Value = '  -0.37%  '

# Row 3
$ws.Range("D3").Value = '3.669.82'
$ws.Range("E3").Value = '  +2.48%  '

# Row 4
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$ws.Range("D5").Value = '''239.69'
$ws.Range("E5").Value = '  -1.03%  '

# Row 6
$ws.Range("D6").Value = '''1.88'
$ws.Range("E6").Value = '  +10.18%  '

# Row 7
$ws.Range("D7").Value = '''655.78'
$ws.Range("E7").Value = '  -0.40%  '

# Row 8
$ws.Range("D8").Value = '''0.424'
$ws.Range("E8").Value = '  -0.65%  '

# Row 9
$ws.Range("D9").Value = '''1.09'
$ws.Range("E9").Value = '  +3.75%  '

# Row 10
$ws.Range("E10").Value = '  -0.01%  '

# Row 11
$ws.Range("D11").Value = '3.668.35'
$ws.Range("E11").Value = '  +2.55%  '

# Row 12
$ws.Range("D12").Value = '''45.51'
$ws.Range("E12").Value = '  +2.32%  '

# Row 13
$ws.Range("E13").Value = '  +1.01%  '

# Row 14
$ws.Range("D14").Value = '''6.84'
$ws.Range("E14").Value = '  +6.41%  '

# Row 15
$ws.Range("D15").Value = '4.354.79'
$ws.Range("E15").Value = '  +2.51%  '

# Row 16
$ws.Range("E16").Value = '  +3.23%  '

# Row 17
$ws.Range("D17").Value = '96.613.62'
$ws.Range("E17").Value = '  -0.64%  '

# Row 18
$ws.Range("D18").Value = '''8.93'
$ws.Range("E18").Value = '  +3.19%  '

# Row 19
$ws.Range("D19").Value = '3.672.46'
$ws.Range("E19").Value = '  +2.69%  '

# Row 20
$ws.Range("D20").Value = '''18.90'
$ws.Range("E20").Value = '  +4.83%  '

# Row 21
$ws.Range("D21").Value = '''12.78'
$ws.Range("E21").Value = '  +0.53%  '

# Row 22
$ws.Range("D22").Value = '''0.531'
$ws.Range("E22").Value = '  +0.71%  '

# Row 23
$ws.Range("D23").Value = '''532.70'
$ws.Range("E23").Value = '  +3.18%  '

# Row 24
$ws.Range("E24").Value = '  +0.11%  '

# Row 25
$ws.Range("E25").Value = '  +5.23%  '

# Row 26
$ws.Range("E26").Value = '  -0.26%  '

# Row 27
$ws.Range("D27").Value = '''102.37'
$ws.Range("E27").Value = '  +0.92%  '

# Row 28
$ws.Range("D28").Value = '''13.52'
$ws.Range("E28").Value = '  +3.36%  '

# Row 29
$ws.Range("D29").Value = '''0.167'
$ws.Range("E29").Value = '  +1.85%  '

# Row 30
$ws.Range("D30").Value = '''12.46'
$ws.Range("E30").Value = '  +4.60%  '

# Row 31
$ws.Range("D31").Value = '''3.05'
$ws.Range("E31").Value = '  +1.59%  '

# Row 32
$ws.Range("D32").Value = '''1.00'
$ws.Range("E32").Value = '  +0.31%  '

# Row 33
$ws.Range("D33").Value = '''1.89'
$ws.Range("E33").Value = '  +14.73%  '

# Row 34
$ws.Range("D34").Value = '''0.186'
$ws.Range("E34").Value = '  +0.87%  '

# Row 35
$ws.Range("D35").Value = '''0.999'
$ws.Range("E35").Value = '  +0.17%  '

# Row 36
$ws.Range("D36").Value = '''32.74'
$ws.Range("E36").Value = '  +3.16%  '

# Row 37
$ws.Range("D37").Value = '''653.94'
$ws.Range("E37").Value = '  +5.58%  '

# Row 38
$ws.Range("D38").Value = '''0.599'
$ws.Range("E38").Value = '  +5.34%  '

# Row 39
$ws.Range("E39").Value = '  +0.24%  '

# Row 40
$ws.Range("E40").Value = '  +4.99%  '

# Row 41
$ws.Range("D41").Value = '''6.83'
$ws.Range("E41").Value = '  +13.88%  '

# Row 42
$ws.Range("D42").Value = '''2.00'
$ws.Range("E42").Value = '  +2.51%  '

# Row 43
$ws.Range("E43").Value = '  +4.26%  '

# Row 44
$ws.Range("D44").Value = '''38.86'
$ws.Range("E44").Value = '  +17.93%  '

# Row 46
$ws.Range("D46").Value = '''0.0459'
$ws.Range("E46").Value = '  +4.44%  '

# Row 47
$ws.Range("D47").Value = '''0.442'
$ws.Range("E47").Value = '  +9.04%  '

# Row 48
$ws.Range("D48").Value = '''2.34'
$ws.Range("E48").Value = '  +1.54%  '

# Row 49
$ws.Range("E49").Value = '  +0.08%  '

# Row 50
$ws.Range("D50").Value = '''8.74'
$ws.Range("E50").Value = '  +2.87%  '

# Row 51
$ws.Range("D51").Value = '''3.65'
$ws.Range("E51").Value = '  +4.35%  '
